$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 122.253015
$ws.Range("H2").Value = 366.759045
$ws.Range("I2").Value = 0.1988639364328829
$ws.Range("J2").Value = 0.1988639364328829
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 12096.53533337286
$ws.Range("R2").Value = 108868.8180003557
$ws.Range("S2").Value = 0.04172184662223347
$ws.Range("T2").Value = 0.04172184662223347

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 122.253015
$ws.Range("H3").Value = 366.759045
$ws.Range("I3").Value = 0.1988639364328829
$ws.Range("J3").Value = 0.1988639364328829
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 19928.00377405054
$ws.Range("R3").Value = 179352.0339664548
$ws.Range("S3").Value = 0.06873316152389551
$ws.Range("T3").Value = 0.06873316152389553

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 122.253015
$ws.Range("H4").Value = 366.759045
$ws.Range("I4").Value = 0.1988639364328829
$ws.Range("J4").Value = 0.1988639364328829
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 7994.871209261602
$ws.Range("R4").Value = 71953.84088335441
$ws.Range("S4").Value = 0.02757490315735857
$ws.Range("T4").Value = 0.02757490315735857

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 122.253015
$ws.Range("H5").Value = 366.759045
$ws.Range("I5").Value = 0.1988639364328829
$ws.Range("J5").Value = 0.1988639364328829
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 17637.78437498195
$ws.Range("R5").Value = 158740.0593748376
$ws.Range("S5").Value = 0.06083402512939531
$ws.Range("T5").Value = 0.06083402512939531

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 132.5447616666667
$ws.Range("H6").Value = 397.634285
$ws.Range("I6").Value = 0.2156050961899926
$ws.Range("J6").Value = 0.2156050961899926
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 13114.86995027745
$ws.Range("R6").Value = 118033.829552497
$ws.Range("S6").Value = 0.04523415816646449
$ws.Range("T6").Value = 0.04523415816646449

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 132.5447616666667
$ws.Range("H7").Value = 397.634285
$ws.Range("I7").Value = 0.2156050961899926
$ws.Range("J7").Value = 0.2156050961899926
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 21605.62265661883
$ws.Range("R7").Value = 194450.6039095695
$ws.Range("S7").Value = 0.07451939334814142
$ws.Range("T7").Value = 0.07451939334814144

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 132.5447616666667
$ws.Range("H8").Value = 397.634285
$ws.Range("I8").Value = 0.2156050961899926
$ws.Range("J8").Value = 0.2156050961899926
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 8667.911372061246
$ws.Range("R8").Value = 78011.2023485512
$ws.Range("S8").Value = 0.02989626854579827
$ws.Range("T8").Value = 0.02989626854579828

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 132.5447616666667
$ws.Range("H9").Value = 397.634285
$ws.Range("I9").Value = 0.2156050961899926
$ws.Range("J9").Value = 0.2156050961899926
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 19122.60344916679
$ws.Range("R9").Value = 172103.4310425011
$ws.Range("S9").Value = 0.06595527612958839
$ws.Range("T9").Value = 0.0659552761295884

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 320.0894206666666
$ws.Range("H10").Value = 960.2682619999999
$ws.Range("I10").Value = 0.5206762565675317
$ws.Range("J10").Value = 0.5206762565675317
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 31671.79956202456
$ws.Range("R10").Value = 285046.196058221
$ws.Range("S10").Value = 0.1092383833188427
$ws.Range("T10").Value = 0.1092383833188428

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 320.0894206666666
$ws.Range("H11").Value = 960.2682619999999
$ws.Range("I11").Value = 0.5206762565675317
$ws.Range("J11").Value = 0.5206762565675317
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 52176.57153959746
$ws.Range("R11").Value = 469589.1438563772
$ws.Range("S11").Value = 0.1799608611106412
$ws.Range("T11").Value = 0.1799608611106412

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 320.0894206666666
$ws.Range("H12").Value = 960.2682619999999
$ws.Range("I12").Value = 0.5206762565675317
$ws.Range("J12").Value = 0.5206762565675317
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 20932.60189678887
$ws.Range("R12").Value = 188393.4170710999
$ws.Range("S12").Value = 0.07219809488198177
$ws.Range("T12").Value = 0.07219809488198178

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 320.0894206666666
$ws.Range("H13").Value = 960.2682619999999
$ws.Range("I13").Value = 0.5206762565675317
$ws.Range("J13").Value = 0.5206762565675317
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 46180.19590299312
$ws.Range("R13").Value = 415621.7631269382
$ws.Range("S13").Value = 0.1592789172560659
$ws.Range("T13").Value = 0.1592789172560659

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 39.86989333333333
$ws.Range("H14").Value = 119.60968
$ws.Range("I14").Value = 0.06485471080959287
$ws.Range("J14").Value = 0.06485471080959287
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 3944.995331562773
$ws.Range("R14").Value = 35504.95798406496
$ws.Range("S14").Value = 0.01360658119145889
$ws.Range("T14").Value = 0.01360658119145889

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 39.86989333333333
$ws.Range("H15").Value = 119.60968
$ws.Range("I15").Value = 0.06485471080959287
$ws.Range("J15").Value = 0.06485471080959287
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 6499.041228698195
$ws.Range("R15").Value = 58491.37105828376
$ws.Range("S15").Value = 0.02241567472524489
$ws.Range("T15").Value = 0.0224156747252449

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 39.86989333333333
$ws.Range("H16").Value = 119.60968
$ws.Range("I16").Value = 0.06485471080959287
$ws.Range("J16").Value = 0.06485471080959287
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 2607.335797215289
$ws.Range("R16").Value = 23466.0221749376
$ws.Range("S16").Value = 0.008992894347520855
$ws.Range("T16").Value = 0.008992894347520857

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 39.86989333333333
$ws.Range("H17").Value = 119.60968
$ws.Range("I17").Value = 0.06485471080959287
$ws.Range("J17").Value = 0.06485471080959287
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 5752.141014001689
$ws.Range("R17").Value = 51769.2691260152
$ws.Range("S17").Value = 0.01983956054536823
$ws.Range("T17").Value = 0.01983956054536823
